$wb = $excel.ActiveWorkbook

# --- Productdata sheet: StartingInventories (C) and SetupCosts (E) ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("C2").Value = 0
$wsProd.Range("E2").Value = 192.456
$wsProd.Range("C3").Value = 1819
$wsProd.Range("E3").Value = 69.92999999999999
$wsProd.Range("C4").Value = 469
$wsProd.Range("E4").Value = 17.50333333333333
$wsProd.Range("C5").Value = 0
$wsProd.Range("E5").Value = 24.37516666666667
$wsProd.Range("C6").Value = 919
$wsProd.Range("E6").Value = 34.16333333333333
$wsProd.Range("C7").Value = 289
$wsProd.Range("E7").Value = 10.272
$wsProd.Range("C8").Value = 109
$wsProd.Range("E8").Value = 3.409333333333333
$wsProd.Range("C9").Value = 1267
$wsProd.Range("E9").Value = 50.49566666666666
$wsProd.Range("C10").Value = 911
$wsProd.Range("E10").Value = 35.99999999999999
$wsProd.Range("C11").Value = 2719
$wsProd.Range("E11").Value = 107.445
$wsProd.Range("C12").Value = 649
$wsProd.Range("E12").Value = 25.235
$wsProd.Range("C13").Value = 0
$wsProd.Range("E13").Value = 184.338
$wsProd.Range("C14").Value = 0
$wsProd.Range("E14").Value = 68.41666666666666
$wsProd.Range("C15").Value = 0
$wsProd.Range("E15").Value = 16.87416666666667
$wsProd.Range("C16").Value = 0
$wsProd.Range("E16").Value = 23.37766666666666
$wsProd.Range("C17").Value = 0
$wsProd.Range("E17").Value = 33.49499999999999
$wsProd.Range("C18").Value = 0
$wsProd.Range("E18").Value = 10.071
$wsProd.Range("C19").Value = 0
$wsProd.Range("E19").Value = 3.343333333333333
$wsProd.Range("C20").Value = 0
$wsProd.Range("E20").Value = 43.83333333333334
$wsProd.Range("C21").Value = 0
$wsProd.Range("E21").Value = 46.53333333333333
$wsProd.Range("C22").Value = 0
$wsProd.Range("E22").Value = 57.59999999999999
$wsProd.Range("C23").Value = 0
$wsProd.Range("E23").Value = 177.2333333333333

# --- Capacity sheet: column B ---
$wsCap = $wb.Worksheets.Item("Capacity")
$wsCap.Range("B2").Value = 16200
$wsCap.Range("B3").Value = 6000
$wsCap.Range("B4").Value = 7500
$wsCap.Range("B5").Value = 6300
$wsCap.Range("B6").Value = 15000
$wsCap.Range("B7").Value = 900
$wsCap.Range("B8").Value = 1500
$wsCap.Range("B9").Value = 12600
$wsCap.Range("B10").Value = 15000
$wsCap.Range("B11").Value = 18000
$wsCap.Range("B12").Value = 8400
$wsCap.Range("B13").Value = 64800
$wsCap.Range("B14").Value = 12000
$wsCap.Range("B15").Value = 4500
$wsCap.Range("B16").Value = 2100
$wsCap.Range("B17").Value = 12000
$wsCap.Range("B18").Value = 900
$wsCap.Range("B19").Value = 900
$wsCap.Range("B20").Value = 30000
$wsCap.Range("B21").Value = 90000
$wsCap.Range("B22").Value = 150000
$wsCap.Range("B23").Value = 150000

# --- ProcessingTime sheet: scattered cells ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")
$wsProc.Range("B2").Value = 1
$wsProc.Range("D4").Value = 5
$wsProc.Range("F6").Value = 5
$wsProc.Range("G7").Value = 1
$wsProc.Range("H8").Value = 5
$wsProc.Range("I9").Value = 3
$wsProc.Range("K11").Value = 2
$wsProc.Range("N14").Value = 2
$wsProc.Range("O15").Value = 3
$wsProc.Range("P16").Value = 1
$wsProc.Range("Q17").Value = 4
$wsProc.Range("R18").Value = 1
$wsProc.Range("S19").Value = 3
$wsProc.Range("T20").Value = 1
$wsProc.Range("V22").Value = 5
